$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217 (shifts existing rows 217:300 down to 218:301,
# extends the used range to A1:R301, just like a manual Excel row insert).
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row 217 with the new weekly price record.
$ws.Cells.Item(217, 1).Value = 8
$ws.Cells.Item(217, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(217, 3).Value = "Coquimbo"
$ws.Cells.Item(217, 4).Value = 44784
$ws.Cells.Item(217, 5).Value = 4
$ws.Cells.Item(217, 6).Value = 100112012
$ws.Cells.Item(217, 7).Value = "Espinaca"
$ws.Cells.Item(217, 8).Value = "Sin especificar"
$ws.Cells.Item(217, 9).Value = "Primera"
$ws.Cells.Item(217, 10).Value = 2400
$ws.Cells.Item(217, 11).Value = 500
$ws.Cells.Item(217, 12).Value = 600
$ws.Cells.Item(217, 13).Value = 550
$ws.Cells.Item(217, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(217, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(217, 16).Value = 1100
$ws.Cells.Item(217, 17).Value = 0.5
$ws.Cells.Item(217, 18).Value = "Hortaliza"
